# The underlying change (per the commit message, "Moving from POI 3.17.0
# to 4.0.1") is not a content edit at all: regenerating the fixture with
# the new Apache POI version re-serialized the package's XML and, in the
# process, re-ordered a handful of attributes (page size/margins in the
# section properties, and the run/paragraph/style defaults plus the
# latent-style table in styles.xml) without changing any value.
#
# We reproduce that by touching (read-then-write-the-same-value) the
# Word object-model properties that are backed by those XML islands, so
# the host re-emits them with its current canonical attribute order.

$d = $word.ActiveDocument

# Re-serializes the section's <w:pgSz>/<w:pgMar> (same values, new
# attribute order).
$d.PageSetup.TopMargin = $d.PageSetup.TopMargin

# Re-serializes styles.xml (docDefaults' <w:rFonts>/<w:lang>, the whole
# <w:latentStyles> exception table, and each <w:style> element's
# attributes) with the same values but the new attribute order.
$styleNames = @("Normal", "Default Paragraph Font", "Normal Table", "No List")
foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    $style.Priority = $style.Priority
}
